$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: drop all <w:lang w:val="fr-FR"/> and merge the trailing
#    " " + "integration" runs (dropping their proofErr spell-check wraps)
#    into a single " integration" run.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$titleXml = '<w:p w:rsidR="001E4365" w:rsidRPr="001E4365" w:rsidRDefault="001E4365" w:rsidP="001E4365"><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/></w:pBdr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr><w:r w:rsidRPr="001E4365"><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve">ALTUI / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="001E4365"><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Thingspeak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="001E4365"><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve"> integration</w:t></w:r></w:p>'
$p1.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) The empty paragraph right after the title loses its pPr/rPr (including
#    its <w:lang w:val="fr-FR"/>) entirely, becoming a bare <w:p/>.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML('<w:p/>')

# ---------------------------------------------------------------------------
# 3) "Step 1 : Preparing a channel in Thingspeak" heading: merge the
#    standalone " " run into the following "Preparing a channel in " run.
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$step1Xml = '<w:p w:rsidR="00E718C0" w:rsidRPr="001E4365" w:rsidRDefault="00E718C0"><w:pPr><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00E718C0"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">Step </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00E718C0"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>1 :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00E718C0"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> Preparing a channel in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00E718C0"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Thingspeak</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p8.Range.InsertXML($step1Xml)

# ---------------------------------------------------------------------------
# 4) Tail section: the "...icon again.  It will show the graph" paragraph
#    gains new trailing sentences (and loses its bookmark), two new
#    paragraphs are inserted, and the picture (now living with the bookmark
#    in the last of those new paragraphs) is cropped/resized and loses its
#    <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$p42 = $d.Paragraphs.Item(42)
$p43 = $d.Paragraphs.Item(43)
$tailRng = $d.Range($p42.Range.Start, $p43.Range.End)
$tailXml = '<w:p w:rsidR="00E718C0" w:rsidRDefault="001E4365" w:rsidP="001E4365"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Once variable start to send data to </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Thingspeak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  (</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> at each change )  , you can go back to the Device Variable and open the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>barGraph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> icon again.  It will show the graph</w:t></w:r><w:r><w:t xml:space="preserve">.  If you leave that page / graph </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>open ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> you will see it adding the points dynamically into it. Quite fun</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:noProof/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="621FD397" wp14:editId="7A052BF5"><wp:extent cx="5937979" cy="3076575"/><wp:effectExtent l="0" t="0" r="5715" b="0"/><wp:docPr id="12" name="Picture 12"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill rotWithShape="1"><a:blip r:embed="rId16"/><a:srcRect b="34483"/><a:stretch/></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3079487"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:ln><a:noFill/></a:ln><a:extLst><a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}"><a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/></a:ext></a:extLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
$tailRng.InsertXML($tailXml)

Write-Output "done"
